$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 43. This shifts the existing rows 43-128 down to 44-129,
# exactly like the source diff shows (old row 43 data now lives at row 44, etc.)
$ws.Rows("43:43").Insert()

# The newly inserted row 43 is mostly blank except for the style carried on column D.
# Populate it by duplicating the row that was just pushed down (now row 44), which
# matches the unchanged columns (A, B, C, E-J, L, Q, T) of the new record in the diff.
$ws.Rows("44:44").Copy()
$ws.Rows("43:43").PasteSpecial()

# Now overwrite the specific cells of the new row 43 with the new record's data.
$ws.Range("D43").Value = 44987
$ws.Range("K43").Value = "Valencia"
$ws.Range("M43").Value = 330
$ws.Range("N43").Value = 1100
$ws.Range("O43").Value = 1150
$ws.Range("P43").Value = 1124
$ws.Range("R43").Value = "Región Metropolitana"
$ws.Range("S43").Value = 1124
